$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# Simple single-value cell updates (rows are 1-based in the Word table)
$tbl.Cell(1,1).Range.Text  = "0M"
$tbl.Cell(2,1).Range.Text  = "0M"
$tbl.Cell(3,1).Range.Text  = "0M"
$tbl.Cell(4,1).Range.Text  = "84"
$tbl.Cell(6,1).Range.Text  = "0.00050"
$tbl.Cell(7,1).Range.Text  = "0.00021"
$tbl.Cell(8,1).Range.Text  = "0.00006"
$tbl.Cell(9,1).Range.Text  = "0.00029"
$tbl.Cell(10,1).Range.Text = "0.00037"
$tbl.Cell(11,1).Range.Text = "0.00041"
$tbl.Cell(12,1).Range.Text = "0.01751"

# Rows that previously held a tab-separated list of values now collapse to a
# single value matching the first column's updated figure.
$tbl.Cell(44,1).Range.Text = "99.97"
$tbl.Cell(45,1).Range.Text = "0.02"
$tbl.Cell(46,1).Range.Text = "62"
